$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column for all existing data rows (2-331)
#    from 45192 (2023-09-23) to 45202 (2023-10-03).
$ws.Range("C2:C331").Value = 45202

# 2. Row 331 picks up an explicit custom height (matches every other data row).
$ws.Rows.Item(331).RowHeight = 15

# 3. Append the new record as row 332.
$ws.Range("A332").Value = "A 45787-2023"

$ws.Range("B332").Value = 45195
$ws.Range("B332").NumberFormat = "YYYY-MM-DD"

$ws.Range("C332").Value = 45202
$ws.Range("C332").NumberFormat = "YYYY-MM-DD"

$ws.Range("D332").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E332").Value = "TRANÅS"

$ws.Range("G332").Value = 9.300000000000001
$ws.Range("H332").Value = 0
$ws.Range("I332").Value = 0
$ws.Range("J332").Value = 0
$ws.Range("K332").Value = 0
$ws.Range("L332").Value = 0
$ws.Range("M332").Value = 0
$ws.Range("N332").Value = 0
$ws.Range("O332").Value = 0
$ws.Range("P332").Value = 0
$ws.Range("Q332").Value = 0

$ws.Range("R332").WrapText = $true
